$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '61.392.59'
Set-TextCell $ws 'E2' '  -3.91%  '
# Row 3
Set-TextCell $ws 'D3' '2.988.43'
Set-TextCell $ws 'E3' '  -3.46%  '
# Row 4
Set-TextCell $ws 'E4' '  +0.06%  '
# Row 5
Set-TextCell $ws 'D5' '535.31'
Set-TextCell $ws 'E5' '  -0.91%  '
# Row 6
Set-TextCell $ws 'D6' '132.77'
Set-TextCell $ws 'E6' '  -3.49%  '
# Row 7
Set-TextCell $ws 'E7' '  +0.03%  '
# Row 8
Set-TextCell $ws 'D8' '2.984.04'
Set-TextCell $ws 'E8' '  -3.39%  '
# Row 9
Set-TextCell $ws 'E9' '  -0.59%  '
# Row 10
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D10' '6.13'
Set-TextCell $ws 'E10' '  -3.66%  '
# Row 11
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws 'D11' '0.148'
Set-TextCell $ws 'E11' '  -6.19%  '
# Row 12
Set-TextCell $ws 'D12' '0.444'
Set-TextCell $ws 'E12' '  -3.76%  '
# Row 13
Set-TextCell $ws 'E13' '  -3.50%  '
# Row 14
Set-TextCell $ws 'D14' '33.59'
Set-TextCell $ws 'E14' '  -4.12%  '
# Row 15
Set-TextCell $ws 'D15' '3.481.85'
Set-TextCell $ws 'E15' '  -3.04%  '
# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D16' '0.110'
Set-TextCell $ws 'E16' '  -1.79%  '
# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D17' '61.448.13'
Set-TextCell $ws 'E17' '  -3.73%  '
# Row 18
Set-TextCell $ws 'D18' '3.001.88'
Set-TextCell $ws 'E18' '  -2.97%  '
# Row 19
Set-TextCell $ws 'D19' '6.59'
Set-TextCell $ws 'E19' '  -2.21%  '
# Row 20
Set-TextCell $ws 'D20' '465.25'
Set-TextCell $ws 'E20' '  -5.15%  '
# Row 21
Set-TextCell $ws 'D21' '13.16'
Set-TextCell $ws 'E21' '  -2.94%  '
# Row 22
Set-TextCell $ws 'D22' '0.669'
Set-TextCell $ws 'E22' '  -5.14%  '
# Row 23
Set-TextCell $ws 'E23' '  -4.67%  '
# Row 24
Set-TextCell $ws 'D24' '80.38'
Set-TextCell $ws 'E24' '  +0.31%  '
# Row 25
Set-TextCell $ws 'D25' '11.89'
Set-TextCell $ws 'E25' '  -3.26%  '
# Row 26
Set-TextCell $ws 'D26' '0.998'
Set-TextCell $ws 'E26' '  -0.38%  '
# Row 27
Set-TextCell $ws 'E27' '  -2.48%  '
# Row 28
Set-TextCell $ws 'D28' '7.67'
Set-TextCell $ws 'E28' '  -7.97%  '
# Row 29
Set-TextCell $ws 'E29' '  +0.36%  '
# Row 30
Set-TextCell $ws 'D30' '1.15'
Set-TextCell $ws 'E30' '  +2.66%  '
# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D31' '1.87'
Set-TextCell $ws 'E31' '  -2.50%  '
# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D32' '25.47'
Set-TextCell $ws 'E32' '  -3.25%  '
# Row 33
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D33' '55.26'
Set-TextCell $ws 'E33' '  -3.32%  '
# Row 34
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D34' '5.42'
Set-TextCell $ws 'E34' '  -1.60%  '
# Row 35
$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D35' '2.26'
Set-TextCell $ws 'E35' '  -6.56%  '
# Row 36
Set-TextCell $ws 'D36' '5.86'
Set-TextCell $ws 'E36' '  -4.03%  '
# Row 37
Set-TextCell $ws 'D37' '449.80'
Set-TextCell $ws 'E37' '  -10.12%  '
# Row 38
Set-TextCell $ws 'D38' '3.158.16'
Set-TextCell $ws 'E38' '  -4.87%  '
# Row 39
Set-TextCell $ws 'D39' '0.0783'
Set-TextCell $ws 'E39' '  -2.60%  '
# Row 40
Set-TextCell $ws 'D40' '0.0382'
Set-TextCell $ws 'E40' '  -4.72%  '
# Row 41
Set-TextCell $ws 'E41' '  +0.01%  '
# Row 42
Set-TextCell $ws 'D42' '8.06'
Set-TextCell $ws 'E42' '  -1.69%  '
# Row 43
Set-TextCell $ws 'E43' '  -10.18%  '
# Row 44
Set-TextCell $ws 'D44' '26.30'
Set-TextCell $ws 'E44' '  +4.03%  '
# Row 45
Set-TextCell $ws 'E45' '  +0.09%  '
# Row 46
Set-TextCell $ws 'D46' '0.241'
Set-TextCell $ws 'E46' '  -7.70%  '
# Row 47
Set-TextCell $ws 'E47' '  -6.64%  '
# Row 48
Set-TextCell $ws 'D48' '118.07'
Set-TextCell $ws 'E48' '  -3.53%  '
# Row 49
Set-TextCell $ws 'E49' '  -2.19%  '
# Row 50
Set-TextCell $ws 'D50' '0.0₃0491'
Set-TextCell $ws 'E50' '  -9.85%  '
# Row 51
Set-TextCell $ws 'E51' '  +6.26%  '
